$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7723323445918061
$ws.Range("C2").Value = 0.03445658242777255
$ws.Range("D2").Value = 0.7420760148450991
$ws.Range("E2").Value = 0.7470135935518118
$ws.Range("F2").Value = 0.7504261789608688
$ws.Range("G2").Value = 0.815963055533745
$ws.Range("H2").Value = 0.8285271493034808
$ws.Range("I2").Value = 12.16421095134998
$ws.Range("J2").Value = 141.3751908297425
$ws.Range("K2").Value = 2.085189552194268
$ws.Range("L2").Value = 67.79968309402466

$ws.Range("B3").Value = 1.029607109549036
$ws.Range("C3").Value = 0.0432930973429417
$ws.Range("D3").Value = 0.9924916943952069
$ws.Range("E3").Value = 0.9978198610884508
$ws.Range("F3").Value = 1.001742467856205
$ws.Range("G3").Value = 1.085878115265738
$ws.Range("H3").Value = 1.098446944070902
$ws.Range("I3").Value = 13.77736584698096
$ws.Range("J3").Value = 734.7185699630767
$ws.Range("K3").Value = 10.83660773080854
$ws.Range("L3").Value = 67.79968309402466
